{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\" \u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u5feb\u901f\u7b80\u5355\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertAfter(\" \u4f7f\u7528git\u521b\u5efa\u5206\u652f\u7b80\u5355\u5feb\u901f\u7b80\u5355\")\n"}
